# Update test case EC_012 (row 13) from Pass -> Fail: the search worked but the
# first link in the results didn't open, so expected/actual result text and the
# status need to be corrected.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I13").Value = "Fail"
$ws.Range("H13").Value = "Chrome browser performed the search, clicked the result, but first link didn't open "
$ws.Range("G13").Value = "The website should search for “Bangalore” and navigate to the clicked search result page and open first link."

# Restore the view to what it was left at after the edit (scrolled up a bit and
# a different active cell selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("K16").Select()
